# Auto update Excel log
# Appends new PRESENCE_DETECTED rows to the "mmWave" worksheet (rows 46-52)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-01-31", "22:01:13", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:01:16", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:01:26", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:01:37", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:01:47", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:01:58", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:02:08", "22:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 46
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 6))
    $rng.NumberFormat = "@"
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
